# Auto-generated edit script
# Updates column F (热度/热度值 count) values across the four sheets
# as described in the commit: 'Update gh-pages to output generated at 456a3b4'

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1546
$ws.Range("F3").Value = 911
$ws.Range("F4").Value = 472
$ws.Range("F5").Value = 916
$ws.Range("F6").Value = 522
$ws.Range("F7").Value = 7830
$ws.Range("F11").Value = 5674
$ws.Range("F12").Value = 577
$ws.Range("F14").Value = 7902
$ws.Range("F15").Value = 9275
$ws.Range("F17").Value = 923
$ws.Range("F18").Value = 4535
$ws.Range("F19").Value = 687
$ws.Range("F20").Value = 259
$ws.Range("F21").Value = 84
$ws.Range("F22").Value = 292
$ws.Range("F24").Value = 1208
$ws.Range("F25").Value = 125
$ws.Range("F26").Value = 1702
$ws.Range("F27").Value = 737
$ws.Range("F28").Value = 962
$ws.Range("F29").Value = 14
$ws.Range("F30").Value = 1900
$ws.Range("F31").Value = 347
$ws.Range("F32").Value = 2344
$ws.Range("F34").Value = 1497
$ws.Range("F36").Value = 1332
$ws.Range("F38").Value = 804
$ws.Range("F39").Value = 522
$ws.Range("F40").Value = 3015
$ws.Range("F41").Value = 4161
$ws.Range("F42").Value = 198
$ws.Range("F44").Value = 433
$ws.Range("F45").Value = 517
$ws.Range("F48").Value = 181
$ws.Range("F49").Value = 4105

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 28
$ws.Range("F22").Value = 2

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5340

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1546
$ws.Range("F4").Value = 911
$ws.Range("F5").Value = 472
$ws.Range("F6").Value = 916
$ws.Range("F7").Value = 522
$ws.Range("F10").Value = 5674
$ws.Range("F11").Value = 577
$ws.Range("F12").Value = 7902
$ws.Range("F13").Value = 9275
$ws.Range("F14").Value = 28
$ws.Range("F15").Value = 923
$ws.Range("F16").Value = 4535
$ws.Range("F17").Value = 687
$ws.Range("F18").Value = 259
$ws.Range("F19").Value = 84
$ws.Range("F20").Value = 292
$ws.Range("F23").Value = 1208
$ws.Range("F24").Value = 125
$ws.Range("F25").Value = 1702
$ws.Range("F26").Value = 737
$ws.Range("F27").Value = 962
$ws.Range("F28").Value = 14
$ws.Range("F29").Value = 1900
$ws.Range("F30").Value = 347
$ws.Range("F31").Value = 2344
$ws.Range("F36").Value = 804
$ws.Range("F39").Value = 522
$ws.Range("F40").Value = 4161
$ws.Range("F42").Value = 198
$ws.Range("F44").Value = 433
$ws.Range("F45").Value = 517
$ws.Range("F48").Value = 181
$ws.Range("F49").Value = 4105
